$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header for column C: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# Update data rows: the audio-file condition column now holds the
# current training phase identifier instead of a wav file path.
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
